# Updates the cryptocurrency price/volume table with freshly scraped values.
# Generated from the upstream data refresh (GitHub Actions cron job).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number, then the new values for whichever of B/C/D/E changed.
# Price (column D) values that look numeric are apostrophe-prefixed so Excel
# keeps them as text (matching the original inline-string cells) instead of
# coercing them into floating point numbers; the style is then reset to the
# workbook default so no stray "quote prefix" formatting is introduced.
$updates = @(
    @{ Row = 2; D="24.505.32"; E="  -1.29%  " }
    @{ Row = 3; D="1.651.60"; E="  -3.20%  " }
    @{ Row = 4; D="1.002"; E="  -0.42%  " }
    @{ Row = 5; D="313.24"; E="  +0.56%  " }
    @{ Row = 6; D="0.9989"; E="  -0.11%  " }
    @{ Row = 7; D="0.3655"; E="  -2.43%  " }
    @{ Row = 8; D="46.62"; E="  -5.89%  " }
    @{ Row = 9; D="0.3248"; E="  -5.55%  " }
    @{ Row = 10; D="1.125"; E="  -7.02%  " }
    @{ Row = 11; D="0.07029"; E="  -6.61%  " }
    @{ Row = 12; D="0.9987"; E="  -0.35%  " }
    @{ Row = 13; D="5.964"; E="  -5.22%  " }
    @{ Row = 14; D="19.38"; E="  -8.38%  " }
    @{ Row = 15; D="6.615"; E="  -6.41%  " }
    @{ Row = 16; D="1.651.56"; E="  -3.27%  " }
    @{ Row = 17; D="0.00001041"; E="  -7.66%  " }
    @{ Row = 18; D="0.06578"; E="  -1.89%  " }
    @{ Row = 19; D="0.9977"; E="  -0.21%  " }
    @{ Row = 20; D="78.63"; E="  -6.58%  " }
    @{ Row = 21; D="5.934"; E="  -6.97%  " }
    @{ Row = 22; D="15.67"; E="  -9.08%  " }
    @{ Row = 23; D="12.54"; E="  -4.21%  " }
    @{ Row = 24; D="24.485.02"; E="  -1.42%  " }
    @{ Row = 25; D="2.465"; E="  +0.56%  " }
    @{ Row = 26; D="2.323" }
    @{ Row = 27; D="146.53"; E="  -2.37%  " }
    @{ Row = 28; D="18.58"; E="  -8.90%  " }
    @{ Row = 29; D="1.834.03"; E="  -3.25%  " }
    @{ Row = 30; B="BitcoinCash"; C="https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"; D="124.23"; E="  -6.54%  " }
    @{ Row = 31; B="ImmutableX"; C="https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"; D="1.185"; E="  -3.41%  " }
    @{ Row = 32; D="4.058"; E="  -4.58%  " }
    @{ Row = 33; D="5.721"; E="  -15.93%  " }
    @{ Row = 34; D="0.08449"; E="  -3.67%  " }
    @{ Row = 35; D="1.659"; E="  -5.98%  " }
    @{ Row = 36; E="  -11.33%  " }
    @{ Row = 37; D="5.203"; E="  -7.50%  " }
    @{ Row = 38; D="1.266"; E="  -0.77%  " }
    @{ Row = 39; B="VeChain"; C="https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"; D="0.02232"; E="  -7.24%  " }
    @{ Row = 40; B="Hedera"; C="https://coinranking.com/coin/jad286TjB+hedera-hbar"; D="0.06018"; E="  -9.23%  " }
    @{ Row = 41; D="0.2069"; E="  -6.83%  " }
    @{ Row = 42; D="8.112"; E="  -11.13%  " }
    @{ Row = 43; D="0.9984"; E="  -0.18%  " }
    @{ Row = 44; D="0.5893"; E="  -8.37%  " }
    @{ Row = 45; E="  -0.85%  " }
    @{ Row = 46; D="12.61"; E="  -8.48%  " }
    @{ Row = 47; D="0.5616"; E="  -8.32%  " }
    @{ Row = 48; D="123.39"; E="  -4.71%  " }
    @{ Row = 49; D="1.948"; E="  -8.05%  " }
    @{ Row = 50; D="0.06917"; E="  -5.20%  " }
    @{ Row = 51; D="1.188"; E="  -2.08%  " }
)

foreach ($u in $updates) {
    foreach ($col in @("B", "C", "D", "E")) {
        if ($u.ContainsKey($col)) {
            $ref = "$col$($u.Row)"
            $value = $u[$col]
            if ($col -eq "D" -and $value -match "^[0-9]*\.?[0-9]+$") {
                $ws.Range($ref).Value = "'" + $value
                $ws.Range($ref).Style = "Normal"
            } else {
                $ws.Range($ref).Value = $value
            }
        }
    }
}
